$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 47.23036199999999
$ws.Range("H2").Value = 141.691086
$ws.Range("I2").Value = 0.3244251370417807
$ws.Range("J2").Value = 0.3244251370417807
$ws.Range("M2").Value = 1.090291
$ws.Range("N2").Value = 3.270873
$ws.Range("O2").Value = 0.02878663098006733
$ws.Range("P2").Value = 0.02878663098006733
$ws.Range("Q2").Value = 51.49483861534199
$ws.Range("R2").Value = 463.4535475380779
$ws.Range("S2").Value = 0.009339106700679514
$ws.Range("T2").Value = 0.009339106700679516

$ws.Range("G3").Value = 47.23036199999999
$ws.Range("H3").Value = 141.691086
$ws.Range("I3").Value = 0.3244251370417807
$ws.Range("J3").Value = 0.3244251370417807
$ws.Range("O3").Value = 0.7239320554917256
$ws.Range("P3").Value = 0.7239320554917257
$ws.Range("Q3").Value = 1295.00268342732
$ws.Range("R3").Value = 11655.02415084588
$ws.Range("S3").Value = 0.2348617563118411
$ws.Range("T3").Value = 0.2348617563118411

$ws.Range("G4").Value = 47.23036199999999
$ws.Range("H4").Value = 141.691086
$ws.Range("I4").Value = 0.3244251370417807
$ws.Range("J4").Value = 0.3244251370417807
$ws.Range("O4").Value = 0.2472813135282071
$ws.Range("P4").Value = 0.2472813135282071
$ws.Range("Q4").Value = 442.348093514034
$ws.Range("R4").Value = 3981.132841626305
$ws.Range("S4").Value = 0.08022427402926015
$ws.Range("T4").Value = 0.08022427402926015

$ws.Range("I5").Value = 0.4188548944674916
$ws.Range("J5").Value = 0.4188548944674916
$ws.Range("M5").Value = 1.090291
$ws.Range("N5").Value = 3.270873
$ws.Range("O5").Value = 0.02878663098006733
$ws.Range("P5").Value = 0.02878663098006733
$ws.Range("Q5").Value = 66.48333538677633
$ws.Range("R5").Value = 598.3500184809869
$ws.Range("S5").Value = 0.01205742128123072
$ws.Range("T5").Value = 0.01205742128123073

$ws.Range("I6").Value = 0.4188548944674916
$ws.Range("J6").Value = 0.4188548944674916
$ws.Range("O6").Value = 0.7239320554917256
$ws.Range("P6").Value = 0.7239320554917257
$ws.Range("S6").Value = 0.303222484704621
$ws.Range("T6").Value = 0.303222484704621

$ws.Range("I7").Value = 0.4188548944674916
$ws.Range("J7").Value = 0.4188548944674916
$ws.Range("O7").Value = 0.2472813135282071
$ws.Range("P7").Value = 0.2472813135282071
$ws.Range("S7").Value = 0.1035749884816399
$ws.Range("T7").Value = 0.1035749884816399

$ws.Range("I8").Value = 0.2567199684907278
$ws.Range("J8").Value = 0.2567199684907277
$ws.Range("M8").Value = 1.090291
$ws.Range("N8").Value = 3.270873
$ws.Range("O8").Value = 0.02878663098006733
$ws.Range("P8").Value = 0.02878663098006733
$ws.Range("Q8").Value = 40.74824000170867
$ws.Range("R8").Value = 366.734160015378
$ws.Range("S8").Value = 0.007390102998157093
$ws.Range("T8").Value = 0.007390102998157092

$ws.Range("I9").Value = 0.2567199684907278
$ws.Range("J9").Value = 0.2567199684907277
$ws.Range("O9").Value = 0.7239320554917256
$ws.Range("P9").Value = 0.7239320554917257
$ws.Range("S9").Value = 0.1858478144752636
$ws.Range("T9").Value = 0.1858478144752636

$ws.Range("I10").Value = 0.2567199684907278
$ws.Range("J10").Value = 0.2567199684907277
$ws.Range("O10").Value = 0.2472813135282071
$ws.Range("P10").Value = 0.2472813135282071
$ws.Range("S10").Value = 0.06348205101730711
$ws.Range("T10").Value = 0.0634820510173071

